$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 765, shifting existing rows 765-783 down to 767-785.
$ws.Rows("765:766").Insert()

# Populate the two newly inserted rows with this week's data.
$ws.Cells.Item(765, 1).Value = 6
$ws.Cells.Item(765, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(765, 3).Value = 'Metropolitana'
$ws.Cells.Item(765, 4).Value = 45239
$ws.Cells.Item(765, 5).Value = 13
$ws.Cells.Item(765, 6).Value = 'Fruta'
$ws.Cells.Item(765, 7).Value = 100101
$ws.Cells.Item(765, 8).Value = 'Berries'
$ws.Cells.Item(765, 9).Value = 100101001
$ws.Cells.Item(765, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(765, 11).Value = 'Sin especificar'
$ws.Cells.Item(765, 12).Value = 'Especial'
$ws.Cells.Item(765, 13).Value = 480
$ws.Cells.Item(765, 14).Value = 7000
$ws.Cells.Item(765, 15).Value = 7000
$ws.Cells.Item(765, 16).Value = 7000
$ws.Cells.Item(765, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(765, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(765, 19).Value = 4667
$ws.Cells.Item(765, 20).Value = 1.5
$ws.Cells.Item(766, 1).Value = 6
$ws.Cells.Item(766, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(766, 3).Value = 'Metropolitana'
$ws.Cells.Item(766, 4).Value = 45239
$ws.Cells.Item(766, 5).Value = 13
$ws.Cells.Item(766, 6).Value = 'Fruta'
$ws.Cells.Item(766, 7).Value = 100101
$ws.Cells.Item(766, 8).Value = 'Berries'
$ws.Cells.Item(766, 9).Value = 100101001
$ws.Cells.Item(766, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(766, 11).Value = 'Sin especificar'
$ws.Cells.Item(766, 12).Value = 'Especial'
$ws.Cells.Item(766, 13).Value = 350
$ws.Cells.Item(766, 14).Value = 8000
$ws.Cells.Item(766, 15).Value = 8000
$ws.Cells.Item(766, 16).Value = 8000
$ws.Cells.Item(766, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(766, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(766, 19).Value = 4000
$ws.Cells.Item(766, 20).Value = 2
